$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: A1 becomes "line" (moved from its old header slot) ---
$ws.Range("A1").Value = "line"

# --- Row 2 updates (quote-prefixed text triggers the quotePrefix style first) ---
$ws.Range("B2").Value = "1.2.3.4"
$ws.Range("F2").Value = "'="

# --- Row 3 updates ---
$ws.Range("C3").Value = "1.2.3.4"
$ws.Range("F3").Value = "'--"

# --- A1 becomes bold (new font/style created after the quotePrefix style) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.ThemeColor = 1

# --- New row 4 ---
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "1.2.3.4"
$ws.Cells.Item(4,3).Value = "2.2.2.2"
$ws.Cells.Item(4,4).Value = 20
$ws.Cells.Item(4,5).Value = 40
$ws.Cells.Item(4,6).Value = "'itay"

# --- New row 5 ---
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "1.2.3.4"
$ws.Cells.Item(5,3).Value = "4.4.4.4"
$ws.Cells.Item(5,4).Value = 60
$ws.Cells.Item(5,5).Value = 80
$ws.Cells.Item(5,6).Value = "'itay1234"

# --- Selection ---
$null = $ws.Range("D4").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
